$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(6, 7, 8, 10, 12, 13, 14)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = 0
}
